$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.869.97"
$ws.Range("E2").Value = "  -3.36%  "
$ws.Range("D3").Value = "2.568.63"
$ws.Range("E3").Value = "  -1.29%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "508.71"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -2.68%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "144.14"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -6.82%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.06%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.557"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -5.87%  "
$ws.Range("D9").Value = "2.574.06"
$ws.Range("E9").Value = "  -1.31%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "6.23"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -7.11%  "
$ws.Range("E11").Value = "  -3.25%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.332"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -4.42%  "
$ws.Range("E13").Value = "  -0.97%  "
$ws.Range("D14").Value = "3.010.77"
$ws.Range("E14").Value = "  -1.57%  "
$ws.Range("D15").Value = "58.842.19"
$ws.Range("E15").Value = "  -3.43%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "20.63"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -4.99%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.0000135"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -4.58%  "
$ws.Range("D18").Value = "2.565.14"
$ws.Range("E18").Value = "  -1.31%  "
$ws.Range("E19").Value = "  -4.98%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "333.70"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -5.79%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "10.09"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -4.56%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.01%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.97"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -4.11%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "60.03"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.33%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.408"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -4.49%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.14%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.158"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -5.39%  "
$ws.Range("D28").Value = "0.0₃0783"
$ws.Range("E28").Value = "  -7.61%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "6.92"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -6.94%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "5.86"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -6.81%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "149.55"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.90%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "18.59"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -4.19%  "
$ws.Range("E34").Value = "  -3.61%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "3.96"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -5.70%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.901"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -1.99%  "
$ws.Range("E37").Value = "  -8.07%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "35.95"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.46%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.827"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -5.47%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "289.38"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.15%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.39"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -7.96%  "
$ws.Range("E42").Value = "  -7.18%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.608"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -2.48%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0981"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -3.42%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0534"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -4.71%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "18.79"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -4.02%  "
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("E49").Value = "  -4.47%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "4.53"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -8.02%  "
$ws.Range("D51").Value = "1.917.37"
$ws.Range("E51").Value = "  -2.46%  "
